$p = $ppt.ActivePresentation

# Slide 4 = "Purpose - Why?" (the slide that gets the new body-text textboxes).
$s = $p.Slides.Item(4)

# --- TextBox 2: small empty autosize textbox (artifact box left behind by the author) ---
$tb2 = $s.Shapes.AddTextbox(1, 105, 202, 14.545748031496062, 29.081259842519685)
$tb2.TextFrame.WordWrap = 0
$tb2.TextFrame.AutoSize = 1

# --- TextBox 3: the big paragraph of body text ---
$tb3 = $s.Shapes.AddTextbox(1, 83, 191, 804, 290.8125196850394)
$tb3.TextFrame.WordWrap = -1
$tb3.TextFrame.AutoSize = 1

$tr = $tb3.TextFrame.TextRange
$tr.Text = "The purpose of this project is to use various Python and machine learning tools to classify images from the Kaggle competition dataset, Cats and Dogs. "
$tr.InsertAfter("`r`r") | Out-Null
$tr.InsertAfter("For this last phase, the goal was to create a ") | Out-Null
$tr.InsertAfter("PyTorch") | Out-Null
$tr.InsertAfter(" CNN to classify the images and predict the bounding box locations. We were to use Object-Oriented convolutional neural networks in ") | Out-Null
$tr.InsertAfter("PyTorch") | Out-Null
$tr.InsertAfter(".") | Out-Null
$tr.InsertAfter("`r`r") | Out-Null
$tr.InsertAfter("The main aim is to classify images and predict bounding boxes with the best accuracy possible. Also, documentation of work, with graphs and tabular data are included to display results and communicate methods.") | Out-Null
$tr.InsertAfter("`r`r") | Out-Null
$tr.InsertAfter("Other aims include gaining experience with Deep Learning techniques and image ") | Out-Null
$tr.InsertAfter("classification tasks.") | Out-Null

# Re-assert the authored size - the spAutoFit live re-layout (text measurement)
# in this environment doesn't exactly match the original author's rendered
# size, so pin it back to the canonical extent from the source file.
$tb3.Left = 83
$tb3.Top = 191
$tb3.Width = 804
$tb3.Height = 290.8125196850394
